$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D (Price) and E (Volume 1h) hold inline-string text such as "247.55" and
# "1.60%". Pre-formatting the target cells as Text keeps Excel from
# re-interpreting that text as a number/percentage when the new value is
# assigned below (matches the existing inline-string cells in the sheet).
$ws.Range("D2:D23").NumberFormat = "@"
$ws.Range("D25:D28").NumberFormat = "@"
$ws.Range("D40:D45").NumberFormat = "@"
$ws.Range("D47:D50").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E4:E28").NumberFormat = "@"
$ws.Range("E40:E50").NumberFormat = "@"

# Row 2 - BNB
$ws.Range("D2").Value = "247.55"
$ws.Range("E2").Value = "1.60%"

# Row 3 - OKB
$ws.Range("D3").Value = "30.28"

# Row 4 - HuobiToken
$ws.Range("D4").Value = "5.178"
$ws.Range("E4").Value = "0.50%"

# Row 5 - Cronos
$ws.Range("D5").Value = "0.05748"
$ws.Range("E5").Value = "2.32%"

# Row 6 - KuCoinToken
$ws.Range("D6").Value = "6.588"
$ws.Range("E6").Value = "1.53%"

# Row 7 - MXToken
$ws.Range("D7").Value = "0.8611"
$ws.Range("E7").Value = "5.41%"

# Row 8 - FTXToken
$ws.Range("D8").Value = "0.8741"
$ws.Range("E8").Value = "5.06%"

# Row 9 - WazirX
$ws.Range("D9").Value = "0.1366"
$ws.Range("E9").Value = "2.75%"

# Row 10 - MandalaExchangeToken
$ws.Range("D10").Value = "0.06997"
$ws.Range("E10").Value = "1.40%"

# Row 11 - BitrueCoin
$ws.Range("D11").Value = "0.02916"
$ws.Range("E11").Value = "0.86%"

# Row 12 - BitMartToken
$ws.Range("D12").Value = "0.09395"
$ws.Range("E12").Value = "0.22%"

# Row 13 - BitForexToken
$ws.Range("D13").Value = "0.001508"
$ws.Range("E13").Value = "-0.17%"

# Row 14 - CoinExToken
$ws.Range("D14").Value = "0.04125"
$ws.Range("E14").Value = "-9.53%"

# Row 15 - One
$ws.Range("D15").Value = "0.0006015"
$ws.Range("E15").Value = "0.43%"

# Row 16 - TigerCash
$ws.Range("D16").Value = "0.006145"
$ws.Range("E16").Value = "0.49%"

# Row 17 - LEO
$ws.Range("D17").Value = "3.505"
$ws.Range("E17").Value = "-2.92%"

# Row 18 - GateToken
$ws.Range("D18").Value = "3.034"
$ws.Range("E18").Value = "0.36%"

# Row 19 - BTSEToken
$ws.Range("D19").Value = "2.183"
$ws.Range("E19").Value = "-5.36%"

# Row 20 - BitpandaEcosystemToken
$ws.Range("D20").Value = "0.3185"
$ws.Range("E20").Value = "2.35%"

# Row 21 - LiechtensteinCryptoassetsExchange
$ws.Range("D21").Value = "0.03308"
$ws.Range("E21").Value = "6.87%"

# Row 22 - ProBitToken
$ws.Range("D22").Value = "0.1307"
$ws.Range("E22").Value = "1.25%"

# Row 23 - MCDex
$ws.Range("D23").Value = "3.608"
$ws.Range("E23").Value = "-3.55%"

# Row 24 - ZBToken
$ws.Range("E24").Value = "2.73%"

# Row 25 - BitKan
$ws.Range("D25").Value = "0.001213"
$ws.Range("E25").Value = "-0.98%"

# Row 26 - HotbitToken
$ws.Range("D26").Value = "0.004507"
$ws.Range("E26").Value = "0.45%"

# Row 27 - NitroEx
$ws.Range("D27").Value = "0.0001178"
$ws.Range("E27").Value = "20.26%"

# Row 28 - UpBots
$ws.Range("D28").Value = "0.0001390"
$ws.Range("E28").Value = "-0.64%"

# Row 40 - IDEX
$ws.Range("D40").Value = "0.03788"
$ws.Range("E40").Value = "4.17%"

# Rows 41-43 - the generated symbol list reshuffled: BKEXToken/CEJI/KickToken
# rotate down one row (41<-BKEXToken, 42<-CEJI, 43<-KickToken slots become
# KickToken, BKEXToken, CEJI respectively), each with refreshed figures.
# Row 41 - was BKEXToken, now KickToken
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.005682"
$ws.Range("E41").Value = "-6.37%"

# Row 42 - was CEJI, now BKEXToken
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1071"
$ws.Range("E42").Value = "1.94%"

# Row 43 - was KickToken, now CEJI
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.002586"
$ws.Range("E43").Value = "-0.17%"

# Row 44 - LocalTraders
$ws.Range("D44").Value = "0.01013"
$ws.Range("E44").Value = "23.87%"

# Row 45 - CoinLion
$ws.Range("D45").Value = "0.00005109"
$ws.Range("E45").Value = "-3.78%"

# Row 46 - Kangarootoken
$ws.Range("E46").Value = "-0.06%"

# Row 47 - CoinbaseStockToken
$ws.Range("D47").Value = "0.08893"
$ws.Range("E47").Value = "-18.40%"

# Row 48 - BOLO
$ws.Range("D48").Value = "0.002717"
$ws.Range("E48").Value = "2.34%"

# Row 49 - CryptobidCoin
$ws.Range("D49").Value = "0.00002098"
$ws.Range("E49").Value = "-0.06%"

# Row 50 - SpecialPowerGold
$ws.Range("D50").Value = "0.0001998"
$ws.Range("E50").Value = "-0.06%"
